$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Total criame" (M) column: per-row totals = SUM(D:L) ---
# M2 and M3 entered individually (non-shared formulas), M4:M30 filled as one
# operation (creates a shared formula group), matching how a user would type
# the first two rows and then fill down the rest.
$ws.Range("M2").Formula = "=SUM(D2:L2)"
$ws.Range("M3").Formula = "=SUM(D3:L3)"
$ws.Range("M4:M30").Formula = "=SUM(D4:L4)"

# --- Give A1 a bottom border (kept unstyled otherwise) ---
$ws.Range("A1").Borders.Item(9).LineStyle = 1

# --- Header for the new column, copying B1's look (font/alignment/border)
#     then trimming the border down to bottom-only ---
$ws.Range("B1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Borders.LineStyle = -4142
$ws.Range("M1").Borders.Item(9).LineStyle = 1
$ws.Range("M1").Value = "Total criame"

# --- Column sizing: widen E (manually resized by the user) and size the
#     new M column ---
$ws.Columns.Item(5).ColumnWidth = 17.285714285714285
$ws.Columns.Item(13).ColumnWidth = 10

# --- Row 1 got slightly taller after the edits ---
$ws.Rows.Item(1).RowHeight = 36.6

# --- Selection moved ---
[void]$ws.Range("F17").Select()
